$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet2")

# Fill in the previously-missing FDM-Travel airport codes (replacing the
# "FINDES IKKE PÅ FDM-TRAVEL" placeholder) for the destinations that now
# have a known code.
$ws.Range("C10").Value = "DUB"   # Dublin (Irland)
$ws.Range("C12").Value = "REK"   # Reykjavik (Island)
$ws.Range("C19").Value = "MLA"   # Valetta (Malta)
$ws.Range("C21").Value = "VIE"   # Wien (Østrig)
$ws.Range("C35").Value = "STO"   # Stockholm (Sverige)
$ws.Range("C48").Value = "BUD"   # Budapest (Ungarn)
$ws.Range("C50").Value = "DXB"   # Dubai (United Arab Emirates)

# Leave the selection on the last-edited cell, matching the saved view state.
$ws.Range("C50").Select()
